# Applies the "Added many more features" revision to the
# "I Pirati del Bounty" slot review: a new, punchier title/meta-description
# plus rewritten "what we like"/"what we don't like" bullet copy.
#
# Word.Find.Execute signature used below:
#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#           MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#           Format, ReplaceWith, Replace)
# Wrap = 1 (wdFindContinue), Replace = 2 (wdReplaceAll)

$d = $word.ActiveDocument

# Title / H1 heading AND the bolded "CTA" run near the end both carry the
# same old string and both change to the same new string, so one
# find-all/replace-all pass over the whole document handles both.
$d.Content.Find.Execute("Play I Pirati del Bounty Free Slot Review", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Play I Pirati del Bounty for Free", 2)

# "What we like" bullets
$d.Content.Find.Execute("Unique and exciting gameplay mechanics", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Unique gameplay mechanics with different types of Wild symbols", 2)

$d.Content.Find.Execute("Variety of Wild symbols increases chances of winning", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Variable number of free spins awarded by the Scatter symbol", 2)

$d.Content.Find.Execute("Scatter symbol offers variable number of free spins", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Exciting bonus features that increase chances of winning", 2)

$d.Content.Find.Execute("Pirate-themed slot with engaging bonus features", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fresh take on the pirate-themed slot genre", 2)

# "What we don't like" bullet
$d.Content.Find.Execute("Sound effects can be excessively cacophonous", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Cacophonous sound design with noisy reel stops", 2)

# Italic summary / meta-description paragraph at the very end
$d.Content.Find.Execute("Discover the pros and cons of playing I Pirati Del Bounty, an exciting pirate-themed slot game with unique gameplay mechanics and bonus features. Play the game for free now.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Read our review of I Pirati del Bounty and play this exciting pirate-themed slot for free.", 2)
